# Applies the "added harvard case classification" recalculation update:
#  - swaps the average_doctor / average_doctor_old column headers (BP1/BQ1)
#  - shifts each row's previous average_doctor value into average_doctor_old (BQ)
#  - writes freshly recalculated statistics into the *_old metric columns and
#    the new average_doctor column (BP) for data rows 4-13
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: average_doctor / average_doctor_old swap ---
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# --- Row 4 ---
$ws.Range("E4").Value = 0.443
$ws.Range("F4").Value = 0.065
$ws.Range("G4").Value = 0.255
$ws.Range("N4").Value = 0.435
$ws.Range("O4").Value = 0.063
$ws.Range("P4").Value = 0.252
$ws.Range("Q4").Value = 0.019
$ws.Range("R4").Value = 0.014
$ws.Range("S4").Value = 0.119
$ws.Range("W4").Value = 0.305
$ws.Range("X4").Value = 0.114
$ws.Range("Y4").Value = 0.338
$ws.Range("AI4").Value = 0.335
$ws.Range("AJ4").Value = 0.08599999999999999
$ws.Range("AK4").Value = 0.293
$ws.Range("AU4").Value = 0.198
$ws.Range("AV4").Value = 0.03
$ws.Range("AW4").Value = 0.173
$ws.Range("BA4").Value = 1.977
$ws.Range("BB4").Value = 0.149
$ws.Range("BC4").Value = 0.385
$ws.Range("BG4").Value = 0.719
$ws.Range("BH4").Value = 0.142
$ws.Range("BI4").Value = 0.377
$ws.Range("BM4").Value = 0.714
$ws.Range("BN4").Value = 0.074
$ws.Range("BO4").Value = 0.272
$ws.Range("BP4").Value = 0.659
$ws.Range("BQ4").Value = 0.697

# --- Row 5 ---
$ws.Range("E5").Value = 0.5659999999999999
$ws.Range("F5").Value = 0.075
$ws.Range("G5").Value = 0.274
$ws.Range("N5").Value = 0.742
$ws.Range("O5").Value = 0.078
$ws.Range("P5").Value = 0.28
$ws.Range("Q5").Value = 0.01
$ws.Range("R5").Value = 0.003
$ws.Range("S5").Value = 0.051
$ws.Range("W5").Value = 0.282
$ws.Range("X5").Value = 0.105
$ws.Range("Y5").Value = 0.324
$ws.Range("AI5").Value = 0.355
$ws.Range("AJ5").Value = 0.092
$ws.Range("AK5").Value = 0.303
$ws.Range("AU5").Value = 0.374
$ws.Range("AV5").Value = 0.093
$ws.Range("AW5").Value = 0.306
$ws.Range("BA5").Value = 1.337
$ws.Range("BB5").Value = 0.08
$ws.Range("BC5").Value = 0.282
$ws.Range("BG5").Value = 0.388
$ws.Range("BH5").Value = 0.049
$ws.Range("BI5").Value = 0.222
$ws.Range("BM5").Value = 0.553
$ws.Range("BN5").Value = 0.061
$ws.Range("BO5").Value = 0.248
$ws.Range("BP5").Value = 0.446
$ws.Range("BQ5").Value = 0.457

# --- Row 6 ---
$ws.Range("E6").Value = 0.497
$ws.Range("N6").Value = 0.548
$ws.Range("Q6").Value = 0.013
$ws.Range("W6").Value = 0.293
$ws.Range("AI6").Value = 0.345
$ws.Range("AU6").Value = 0.259
$ws.Range("BA6").Value = 1.585
$ws.Range("BG6").Value = 0.504
$ws.Range("BM6").Value = 0.623
$ws.Range("BP6").Value = 0.528
$ws.Range("BQ6").Value = 0.548

# --- Row 7 ---
$ws.Range("E7").Value = 0.536
$ws.Range("N7").Value = 0.65
$ws.Range("Q7").Value = 0.011
$ws.Range("W7").Value = 0.286
$ws.Range("AI7").Value = 0.351
$ws.Range("AU7").Value = 0.318
$ws.Range("BA7").Value = 1.425
$ws.Range("BG7").Value = 0.427
$ws.Range("BM7").Value = 0.579
$ws.Range("BP7").Value = 0.475
$ws.Range("BQ7").Value = 0.489

# --- Row 8 ---
$ws.Range("E8").Value = 0.623
$ws.Range("F8").Value = 0.1
$ws.Range("G8").Value = 0.316
$ws.Range("N8").Value = 0.776
$ws.Range("O8").Value = 0.068
$ws.Range("P8").Value = 0.261
$ws.Range("Q8").Value = 0.01
$ws.Range("S8").Value = 0.077
$ws.Range("W8").Value = 0.317
$ws.Range("X8").Value = 0.121
$ws.Range("Y8").Value = 0.348
$ws.Range("AI8").Value = 0.381
$ws.Range("AJ8").Value = 0.126
$ws.Range("AK8").Value = 0.355
$ws.Range("AU8").Value = 0.317
$ws.Range("AV8").Value = 0.08400000000000001
$ws.Range("AW8").Value = 0.29
$ws.Range("BA8").Value = 1.718
$ws.Range("BB8").Value = 0.12
$ws.Range("BC8").Value = 0.347
$ws.Range("BG8").Value = 0.542
$ws.Range("BH8").Value = 0.108
$ws.Range("BI8").Value = 0.328
$ws.Range("BM8").Value = 0.6919999999999999
$ws.Range("BN8").Value = 0.064
$ws.Range("BO8").Value = 0.253
$ws.Range("BP8").Value = 0.573
$ws.Range("BQ8").Value = 0.595

# --- Row 9 ---
$ws.Range("E9").Value = 0.544
$ws.Range("N9").Value = 0.671
$ws.Range("O9").Value = 0.221
$ws.Range("P9").Value = 0.47
$ws.Range("W9").Value = 0.215
$ws.Range("X9").Value = 0.169
$ws.Range("Y9").Value = 0.411
$ws.Range("AI9").Value = 0.291
$ws.Range("AJ9").Value = 0.206
$ws.Range("AK9").Value = 0.454
$ws.Range("BA9").Value = 1.62
$ws.Range("BB9").Value = 0.245
$ws.Range("BC9").Value = 0.495
$ws.Range("BG9").Value = 0.5570000000000001
$ws.Range("BH9").Value = 0.247
$ws.Range("BI9").Value = 0.497
$ws.Range("BM9").Value = 0.633
$ws.Range("BN9").Value = 0.232
$ws.Range("BO9").Value = 0.482
$ws.Range("BP9").Value = 0.54
$ws.Range("BQ9").Value = 0.5570000000000001

# --- Row 10 ---
$ws.Range("E10").Value = 0.696
$ws.Range("F10").Value = 0.212
$ws.Range("G10").Value = 0.46
$ws.Range("N10").Value = 0.886
$ws.Range("O10").Value = 0.101
$ws.Range("P10").Value = 0.318
$ws.Range("W10").Value = 0.392
$ws.Range("X10").Value = 0.238
$ws.Range("Y10").Value = 0.488
$ws.Range("AI10").Value = 0.418
$ws.Range("AJ10").Value = 0.243
$ws.Range("AK10").Value = 0.493
$ws.Range("AU10").Value = 0.316
$ws.Range("AV10").Value = 0.216
$ws.Range("AW10").Value = 0.465
$ws.Range("BA10").Value = 2.05
$ws.Range("BB10").Value = 0.243
$ws.Range("BC10").Value = 0.493
$ws.Range("BG10").Value = 0.62
$ws.Range("BH10").Value = 0.236
$ws.Range("BI10").Value = 0.485
$ws.Range("BM10").Value = 0.848
$ws.Range("BN10").Value = 0.129
$ws.Range("BO10").Value = 0.359
$ws.Range("BP10").Value = 0.6830000000000001
$ws.Range("BQ10").Value = 0.713

# --- Row 11 ---
$ws.Range("E11").Value = 0.734
$ws.Range("F11").Value = 0.195
$ws.Range("G11").Value = 0.442
$ws.Range("N11").Value = 0.899
$ws.Range("O11").Value = 0.091
$ws.Range("P11").Value = 0.302
$ws.Range("W11").Value = 0.392
$ws.Range("X11").Value = 0.238
$ws.Range("Y11").Value = 0.488
$ws.Range("AI11").Value = 0.456
$ws.Range("AJ11").Value = 0.248
$ws.Range("AK11").Value = 0.498
$ws.Range("AU11").Value = 0.43
$ws.Range("AV11").Value = 0.245
$ws.Range("AW11").Value = 0.495
$ws.Range("BA11").Value = 2.05
$ws.Range("BB11").Value = 0.243
$ws.Range("BC11").Value = 0.493
$ws.Range("BG11").Value = 0.62
$ws.Range("BH11").Value = 0.236
$ws.Range("BI11").Value = 0.485
$ws.Range("BM11").Value = 0.848
$ws.Range("BN11").Value = 0.129
$ws.Range("BO11").Value = 0.359
$ws.Range("BP11").Value = 0.6830000000000001
$ws.Range("BQ11").Value = 0.717

# --- Row 12 ---
$ws.Range("E12").Value = 1.466
$ws.Range("F12").Value = 0.835
$ws.Range("G12").Value = 0.914
$ws.Range("N12").Value = 1.417
$ws.Range("O12").Value = 0.743
$ws.Range("P12").Value = 0.862
$ws.Range("W12").Value = 1.613
$ws.Range("X12").Value = 0.5600000000000001
$ws.Range("Y12").Value = 0.748
$ws.Range("AI12").Value = 1.722
$ws.Range("AJ12").Value = 1.312
$ws.Range("AK12").Value = 1.145
$ws.Range("AU12").Value = 2.722
$ws.Range("AV12").Value = 2.812
$ws.Range("AW12").Value = 1.677
$ws.Range("BA12").Value = 3.805
$ws.Range("BB12").Value = 0.45
$ws.Range("BC12").Value = 0.671
$ws.Range("BG12").Value = 1.122
$ws.Range("BH12").Value = 0.148
$ws.Range("BI12").Value = 0.385
$ws.Range("BM12").Value = 1.313
$ws.Range("BN12").Value = 0.335
$ws.Range("BO12").Value = 0.578
$ws.Range("BP12").Value = 1.268
$ws.Range("BQ12").Value = 1.29

# --- Row 13 ---
$ws.Range("E13").Value = 1.53
$ws.Range("F13").Value = 0.525
$ws.Range("G13").Value = 0.724
$ws.Range("N13").Value = 2.059
$ws.Range("O13").Value = 0.97
$ws.Range("P13").Value = 0.985
$ws.Range("W13").Value = 1.009
$ws.Range("X13").Value = 0.18
$ws.Range("Y13").Value = 0.425
$ws.Range("AI13").Value = 1.265
$ws.Range("AJ13").Value = 0.376
$ws.Range("AK13").Value = 0.613
$ws.Range("AU13").Value = 2.263
$ws.Range("AV13").Value = 0.995
$ws.Range("AW13").Value = 0.998
$ws.Range("BA13").Value = 2.356
$ws.Range("BB13").Value = 0.312
$ws.Range("BC13").Value = 0.5590000000000001
$ws.Range("BG13").Value = 0.58
$ws.Range("BH13").Value = 0.053
$ws.Range("BI13").Value = 0.229
$ws.Range("BM13").Value = 0.885
$ws.Range("BN13").Value = 0.235
$ws.Range("BO13").Value = 0.485
$ws.Range("BP13").Value = 0.785
$ws.Range("BQ13").Value = 0.737
